$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cell values are plain text (inline strings) in the source workbook.
# Some "Price" values look numeric (e.g. "325.48"), so we briefly force a
# text number format before assigning them, then restore the default style
# so no visible formatting change is introduced.

$ws.Range("D2").Value = "30.205.70"
$ws.Range("E2").Value = "  +3.28%  "
$ws.Range("D3").Value = "1.896.82"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "325.48"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +3.44%  "
$ws.Range("E7").Value = "  +0.50%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.4014"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +2.68%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.08446"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +0.30%  "
$ws.Range("E10").Value = "  +0.60%  "
$ws.Range("E11").Value = "  +0.39%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "23.19"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +12.39%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "6.432"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +3.04%  "
$ws.Range("D14").Value = "1.893.71"
$ws.Range("E14").Value = "  +0.19%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "7.342"
$ws.Range("D15").Style = "Normal"
$ws.Range("E16").Value = "  -0.32%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "94.81"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +1.98%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.00001110"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +0.51%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.06658"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -1.20%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "18.27"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +2.59%  "
$ws.Range("E21").Value = "  -0.34%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.948"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -0.99%  "
$ws.Range("D23").Value = "30.214.37"
$ws.Range("E23").Value = "  +3.25%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "11.29"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +1.63%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.228"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.35%  "
$ws.Range("D26").Value = "2.111.74"
$ws.Range("E26").Value = "  +0.24%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "21.71"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +4.13%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "161.24"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +1.10%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.374"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -2.12%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "129.06"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +1.27%  "
$ws.Range("E31").Value = "  +3.38%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.1056"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +0.90%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "6.046"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -2.49%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "3.765"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +3.06%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.02488"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +0.51%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.06545"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -0.03%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.2206"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +0.94%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "5.242"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +2.13%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.219"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -0.62%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "11.79"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +4.86%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.6505"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +0.13%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "8.719"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -3.54%  "
$ws.Range("E43").Value = "  +0.57%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.6110"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +1.15%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "13.28"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +1.20%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "3.704"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +0.89%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.057"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +0.55%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.235"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +0.68%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "124.35"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +0.96%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.164"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -0.98%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "78.88"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +2.07%  "
